$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append an "order By ... LIMIT 100" clause to the three Neo4j queries
# (Cases / Samples / Files tabs) stored in column B.

$caseQuery = $ws.Range("B2").Value2
$ws.Range("B2").Value = $caseQuery + "`n order By ss.study_subject_id ASC LIMIT 100"

$sampleQuery = $ws.Range("B3").Value2
$ws.Range("B3").Value = $sampleQuery + "`n order By samp.sample_id ASC LIMIT 100"

$fileQuery = $ws.Range("B4").Value2
$ws.Range("B4").Value = $fileQuery + "`n order By f.file_name ASC LIMIT 100"

# Update selection to match the authored file (active cell moved to B3,
# and the view no longer pins a frozen top-left cell at A4).
$ws.Range("B3").Select()
